$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 8, shifting existing rows 8-20 down to 9-21
$ws.Rows.Item(8).Insert()

# Populate the newly inserted row 8 with the new weekly data point
$ws.Range("A8").Value = 7
$ws.Range("B8").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C8").Value = "Ñuble"
$ws.Range("D8").Value = 44536
$ws.Range("E8").Value = 16
$ws.Range("F8").Value = 100112040
$ws.Range("G8").Value = "Cilantro"
$ws.Range("H8").Value = "Sin especificar"
$ws.Range("I8").Value = "Primera"
$ws.Range("J8").Value = 50
$ws.Range("K8").Value = 2000
$ws.Range("L8").Value = 2000
$ws.Range("M8").Value = 2000
$ws.Range("N8").Value = "$/atado 0,5 a 1 kilo"
$ws.Range("O8").Value = "Provincia de Diguillín"
$ws.Range("P8").Value = 2000
$ws.Range("Q8").Value = 1
$ws.Range("R8").Value = "Hortaliza"
